$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Header date line
Replace-Text "2025-09-14 Sunday" "2025-09-15 Monday"

# First data row of the table: the "76÷8=" cell is removed, the remaining
# values shift left, and a new "27÷8=" cell is appended at the end -- net
# effect is the same 5 cells with new contents in each position.
$t = $d.Tables.Item(1)
$row1 = $t.Rows.Item(1)
$row1.Cells.Item(1).Range.Text = "10÷8="
$row1.Cells.Item(2).Range.Text = "69÷9="
$row1.Cells.Item(3).Range.Text = "55÷6="
$row1.Cells.Item(4).Range.Text = "15÷2="
$row1.Cells.Item(5).Range.Text = "27÷8="

# Remaining problem cells: straightforward value replacements.
Replace-Text "59÷3=" "83÷5="
Replace-Text "43÷7=" "65÷6="
Replace-Text "57÷4=" "46÷2="
Replace-Text "45÷2=" "39÷9="
Replace-Text "44÷9=" "88÷6="

Replace-Text "89÷6=" "75÷8="
Replace-Text "20÷9=" "26÷9="
Replace-Text "49÷3=" "73÷3="
Replace-Text "64÷4=" "15÷3="
Replace-Text "56÷9=" "78÷9="

Replace-Text "18÷6=" "59÷7="
Replace-Text "90÷8=" "21÷9="
Replace-Text "44÷6=" "87÷2="
Replace-Text "53÷2=" "15÷9="
Replace-Text "21÷3=" "91÷9="

Replace-Text "44÷8=" "64÷3="
Replace-Text "83÷9=" "39÷8="
Replace-Text "71÷6=" "98÷2="
Replace-Text "62÷2=" "79÷7="
Replace-Text "67÷2=" "57÷8="
